# Gummy data key update
# - workbook: active tab -> "Respondent Details" (index 4, 0-based)
# - shared strings: "Min" -> "Min_val", "Max" -> "Max_val" (Panel Details header)
# - Panel Details (sheet3): drop unused extra cols, clear the "extreme bound"
#   number format on D28:E37, refresh the Min/Max values for every attribute,
#   move the selection
# - Sample Info / Scrap: move the selection
# - Respondent Details: becomes the active/selected sheet

$wb = $excel.ActiveWorkbook

$wsSample = $wb.Worksheets.Item("Sample Info")
$wsPanelDetails = $wb.Worksheets.Item("Panel Details")
$wsRespondentDetails = $wb.Worksheets.Item("Respondent Details")
$wsScrap = $wb.Worksheets.Item("Scrap")

# ---- Panel Details header: Min/Max -> Min_val/Max_val ----
$wsPanelDetails.Range("D1").Value = "Min_val"
$wsPanelDetails.Range("E1").Value = "Max_val"

# ---- Panel Details: refreshed Min/Max bounds per attribute ----
$wsPanelDetails.Range("E2").Value = 35.700000000000003
$wsPanelDetails.Range("E3").Value = 21
$wsPanelDetails.Range("E4").Value = 34.65
$wsPanelDetails.Range("E5").Value = 40.950000000000003
$wsPanelDetails.Range("D6").Value = 10.591666666666701
$wsPanelDetails.Range("E6").Value = 52.575000000000003
$wsPanelDetails.Range("E7").Value = 81.900000000000006
$wsPanelDetails.Range("E8").Value = 17.850000000000001
$wsPanelDetails.Range("E9").Value = 70.349999999999994
$wsPanelDetails.Range("E10").Value = 62.475000000000001
$wsPanelDetails.Range("D11").Value = 5.125
$wsPanelDetails.Range("E11").Value = 57.375
$wsPanelDetails.Range("D12").Value = 3.5416666666666701
$wsPanelDetails.Range("E12").Value = 46.625
$wsPanelDetails.Range("E13").Value = 31.5
$wsPanelDetails.Range("D14").Value = 19.7
$wsPanelDetails.Range("E14").Value = 55.633333333333297
$wsPanelDetails.Range("D15").Value = 11.3333333333333
$wsPanelDetails.Range("E15").Value = 55.3333333333333
$wsPanelDetails.Range("E16").Value = 15.75
$wsPanelDetails.Range("E17").Value = 43.05
$wsPanelDetails.Range("D18").Value = 7.6333333333333302
$wsPanelDetails.Range("E18").Value = 37.700000000000003
$wsPanelDetails.Range("D19").Value = 2.18333333333333
$wsPanelDetails.Range("E19").Value = 64.150000000000006
$wsPanelDetails.Range("E20").Value = 19.25
$wsPanelDetails.Range("E21").Value = 27.3
$wsPanelDetails.Range("E22").Value = 21.35
$wsPanelDetails.Range("E23").Value = 14.7
$wsPanelDetails.Range("E24").Value = 35
$wsPanelDetails.Range("E25").Value = 13.8588141025641
$wsPanelDetails.Range("E26").Value = 31.7871794871795
$wsPanelDetails.Range("E27").Value = 17.014743589743599
$wsPanelDetails.Range("D28").Value = -460.27955181715299
$wsPanelDetails.Range("E28").Value = 0
$wsPanelDetails.Range("D29").Value = 38.253010022841302
$wsPanelDetails.Range("E29").Value = 21176.385364753201
$wsPanelDetails.Range("D30").Value = 0.45384820957368799
$wsPanelDetails.Range("E30").Value = 0.87879933661402998
$wsPanelDetails.Range("E31").Value = 25513.339510966802
$wsPanelDetails.Range("E32").Value = 35978.039452392499
$wsPanelDetails.Range("D33").Value = -638.81177887917397
$wsPanelDetails.Range("E33").Value = -276.62056806683199
$wsPanelDetails.Range("D34").Value = -379.56949223212803
$wsPanelDetails.Range("E34").Value = -190.16778661728199
$wsPanelDetails.Range("D35").Value = 18.407415045495
$wsPanelDetails.Range("E35").Value = 57.223971025907097
$wsPanelDetails.Range("D36").Value = 67.092896928727399
$wsPanelDetails.Range("E36").Value = 98.062014109751104
$wsPanelDetails.Range("D37").Value = 216.366249903115
$wsPanelDetails.Range("E37").Value = 31267.260024751999

# The D28:E37 block used to carry a special "extreme bound" number format
# (+/-1E98 sentinel values); now that real bounds are used, drop that
# formatting back to the sheet default.
$wsPanelDetails.Range("D28:E37").ClearFormats()

# Unused helper columns (cols 6-9 / F:I) go away; the remaining D:E pair
# gets a uniform width
$wsPanelDetails.Columns("F:I").Delete()
$wsPanelDetails.Columns("D:E").ColumnWidth = 38.140625

# ---- Selections ----
$wsSample.Range("C2:C10").Select()
$wsPanelDetails.Range("A43").Select()
$wsScrap.Range("G1:G1048576").Select()

# Respondent Details becomes the active sheet/selection
$wsRespondentDetails.Activate()
$wsRespondentDetails.Range("E19").Select()
